$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.945.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.983.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.686"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.759"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.627.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.992.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  +5.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.873.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "49.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "71.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "642.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0930"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0487"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.149"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +35.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.903.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
